# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to match the newly generated output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" - rows 6,8,9,10,11,12 in column F
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F6").Value = 139
$wsExhibit.Range("F8").Value = 4915
$wsExhibit.Range("F9").Value = 107
$wsExhibit.Range("F10").Value = 5208
$wsExhibit.Range("F11").Value = 595
$wsExhibit.Range("F12").Value = 1306

# Sheet "全部类型" - rows 6,9,10,11,12,13 in column F
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F6").Value = 139
$wsAll.Range("F9").Value = 4915
$wsAll.Range("F10").Value = 107
$wsAll.Range("F11").Value = 5208
$wsAll.Range("F12").Value = 595
$wsAll.Range("F13").Value = 1306
